$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the IP-address entries for rows 5 and 6 in column B
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()

# Reflect the final selection from the saved file (active cell B6)
$ws.Range("B6").Select()
